$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰" + "`n" + `
    "✅ Dólar paralelo: 68" + "`n" + `
    "" + "`n" + `
    "Binance" + "`n" + `
    "✅ 1000 Bs = 4.94 = 19614.62 pesos" + "`n" + `
    "✅ 19614.62 pesos = 4.92 = 946.37 Bs" + "`n" + `
    "" + "`n" + `
    "Promedio competencia" + "`n" + `
    "✅ Tasa pesos: 20" + "`n" + `
    "✅ Tasa Bs: 20" + "`n" + `
    "✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 202.4
$ws2.Range("O10").Value = 3970
$ws2.Range("N12").Value = 3990
$ws2.Range("O12").Value = 192.51
